$wb = $excel.ActiveWorkbook

# Add the new "Sources" worksheet at the end of the workbook (after the last
# existing sheet) and make it the active sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Sources"

# General Motors 2017 annual report source.
$ws.Range("B2").Value = "General Motors 2017 annual report:"
$ws.Range("B3").Value = "https://www.gm.com/content/dam/gm/en_us/english/Group4/InvestorsPDFDocuments/02-pdfs/10-K.pdf"

# Fuel type by industry in Denmark source.
$ws.Range("B6").Value = "www.statistikbanken.dk/ENE3H"
$ws.Range("B5").Value = "Fuel type by industry in Denmark:"

# Busiest container ports source.
$ws.Range("B8").Value = "Busiest container ports: "
$ws.Range("B9").Value = "https://en.wikipedia.org/wiki/List_of_busiest_container_ports"

# Busiest airports by passenger traffic source.
$ws.Range("B11").Value = "Busiest airports by passenger traffic"
$ws.Range("B12").Value = "https://en.wikipedia.org/wiki/List_of_busiest_airports_by_passenger_traffic"

$ws.Select() | Out-Null
$ws.Range("E24").Select() | Out-Null
